$d = $word.ActiveDocument
Write-Output $d.Paragraphs.Count
$p = $d.Paragraphs.Item(6)
Write-Output "text=[$($p.Range.Text)]"
